# HOTFIX #5: Fix KULIGI (poi_31) - add missing IDs and adjust opening
# hours to 18:00-22:00 for evening scheduling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row (A1:AQ1): bold, centered/top aligned, thin box border ---
$header = $ws.Range("A1:AQ1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous (thin box around each cell)

# --- 2. Fill in the missing POI IDs in column A (rows 2-32 => poi_1..poi_31) ---
for ($row = 2; $row -le 32; $row++) {
    $ws.Cells.Item($row, 1).Value = "poi_" + ($row - 1)
}

# --- 3. Fix the KULIGI (poi_31, row 32) opening hours ---
$ws.Range("F32").Value = "mon:18:00-22:00,tue:18:00-22:00,wed:18:00-22:00,thu:18:00-22:00,fri:18:00-22:00,sat:18:00-22:00,sun:18:00-22:00"
